# The author shifted the "GAME" step list down by one row: a new blank
# row was opened up above "8. call out player..." (old row 28, which held
# "7. once input has reached players # animation that all questions are in")
# so that every entry from old B28:B33 now lives one row lower, at
# B29:B34. Row 36 ("next things start on players.py line 45") is
# untouched - this was a local shift of just that block, not a sheet-wide
# row insert.
#
# Implement it as a read-all-then-write-shifted-then-clear-source
# operation on column B so every value ends up exactly one row below
# where it started, and the old top row is cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$sourceFirstRow = 28
$sourceLastRow  = 33

# Capture the existing values before any writes happen.
$values = @{}
for ($r = $sourceFirstRow; $r -le $sourceLastRow; $r++) {
    $values[$r] = $ws.Range("B$r").Value()
}

# Write each captured value one row lower than where it came from.
for ($r = $sourceFirstRow; $r -le $sourceLastRow; $r++) {
    $ws.Range("B$($r + 1)").Value = $values[$r]
}

# The old top row (28) is now vacated.
$ws.Range("B$sourceFirstRow").ClearContents()

# Leave the view/selection where the author left it.
$ws.Range("C28").Select()
